# RF013 - Gerenciar Niveis das Competencias: bump from v1.2 to v1.3
#
# TC4 and TC5 swap places in the suite: the old TC4 body (7 rows: header,
# description, precondition, table header + 6 step rows) moves down to the
# TC5 slot (with step #2 now describing the "Usuario Nao-Autenticado" actor
# instead of "Lider de Pessoas ... Novo"), while the old TC5 body (a single
# step) moves up into the now-shorter TC4 slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Wipe the whole affected region (old TC4 rows 40-49 body rows 44-49,
#    blank rows 50-51, old TC5 rows 52-56) of both content and formatting
#    so we can rebuild it row-by-row with the right style + text.
# ---------------------------------------------------------------------
$ws.Range("A44:F56").ClearContents()
$ws.Range("A44:F56").ClearFormats()

# Unmerge the two label rows that belonged to the old TC5 block.
$ws.Range("B53:D53").UnMerge()
$ws.Range("B54:F54").UnMerge()

# ---------------------------------------------------------------------
# 2. New TC4 body: just one step row (row 44), reusing the donor style of
#    any "data row" (row 10 in the TC1 block is a clean, unmerged donor).
# ---------------------------------------------------------------------
$ws.Range("A10:F10").Copy()
$ws.Range("A44:F44").PasteSpecial(-4122)

$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Lider de Pessoas clica na opcao 'Novo' para criar um novo Niveis das Competencias"
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = "SYSTEM exibe uma mensagem de erro ao tentar salvar o novo Nivel das Competencias, informando o campo ou a validacao que falhou"
$ws.Range("E44").Value = ""
$ws.Range("F44").Value = ""

# Rows 45-46 become blank separator rows (already cleared above).

# ---------------------------------------------------------------------
# 3. New TC5 body occupies rows 47-56: header / description / precondition
#    / table header / 6 step rows - i.e. exactly what the old TC4 block
#    (rows 40-49) used to look like.
# ---------------------------------------------------------------------

# Row 47: "Test Case ID:" header row -> donor row 6 (TC1 header).
$ws.Range("A6:F6").Copy()
$ws.Range("A47:F47").PasteSpecial(-4122)
$ws.Range("A47").Value = "Test Case ID: "
$ws.Range("B47").Value = "TC5"
$ws.Range("C47").Value = "Priority (low,medium,high: "
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = "Executed by:"
$ws.Range("F47").Value = ""

# Row 48: "Description:" / "Execution Date:" row -> donor row 7.
$ws.Range("A7:F7").Copy()
$ws.Range("A48:F48").PasteSpecial(-4122)
$ws.Range("A48").Value = "Description: "
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = ""
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = "Execution Date: "
$ws.Range("F48").Value = ""

# Row 49: "Precondition:" row -> donor row 8.
$ws.Range("A8:F8").Copy()
$ws.Range("A49:F49").PasteSpecial(-4122)
$ws.Range("A49").Value = "Precondition: "
$ws.Range("B49").Value = "Lider de Pessoas esta autenticado no sistema e  tem permissao para gerenciar Niveis das Competencias."
$ws.Range("C49").Value = ""
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("F49").Value = ""

# Row 50: "#/Steps/Test Data/Expected Results/..." table header -> donor row 9.
$ws.Range("A9:F9").Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)
$ws.Range("A50").Value = "#"
$ws.Range("B50").Value = "Steps"
$ws.Range("C50").Value = "Test Data"
$ws.Range("D50").Value = "Expected Results"
$ws.Range("E50").Value = "Execution Status (pass/fail/blocked)"
$ws.Range("F50").Value = "Actual Result"

# Rows 51-56: the six step rows, donors from the TC1 data rows (10-15).
$ws.Range("A10:F15").Copy()
$ws.Range("A51:F56").PasteSpecial(-4122)

$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("C51").Value = ""
$ws.Range("D51").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("E51").Value = ""
$ws.Range("F51").Value = ""

$ws.Range("A52").Value = 2
$ws.Range("B52").Value = "Usuario Nao-Autenticado acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("C52").Value = ""
$ws.Range("D52").Value = "SYSTEM apresenta o formulario para cadastro e alteracao de Niveis das Competencias"
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = ""

$ws.Range("A53").Value = 3
$ws.Range("B53").Value = "Lider de Pessoas preenche o campo 'Nome' com dados validos"
$ws.Range("C53").Value = ""
$ws.Range("D53").Value = "SYSTEM apresenta o campo 'Nome' preenchido corretamente"
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = ""

$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Lider de Pessoas preenche o campo 'Valor' com um numero inteiro referente a pontuacao a ser obtida no atendimento do respectivo nivel"
$ws.Range("C54").Value = ""
$ws.Range("D54").Value = "SYSTEM apresenta o campo 'Valor' preenchido corretamente"
$ws.Range("E54").Value = ""
$ws.Range("F54").Value = ""

$ws.Range("A55").Value = 5
$ws.Range("B55").Value = "Lider de Pessoas preenche o campo 'Descricao' com a descricao do nivel de competencia"
$ws.Range("C55").Value = ""
$ws.Range("D55").Value = "SYSTEM apresenta o campo 'Descricao' preenchido corretamente"
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = ""

$ws.Range("A56").Value = 6
$ws.Range("B56").Value = "Lider de Pessoas clica na opcao 'Salvar'"
$ws.Range("C56").Value = ""
$ws.Range("D56").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastrados apenas para visualizacao com a opcao 'Ajuda'"
$ws.Range("E56").Value = ""
$ws.Range("F56").Value = ""

# ---------------------------------------------------------------------
# 4. Merges: the description/precondition rows of the (new) TC5 block
#    live at rows 48/49 now.
# ---------------------------------------------------------------------
$ws.Range("B48:D48").Merge()
$ws.Range("B49:F49").Merge()
